$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (column F = "想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9
$ws1.Range("F3").Value = 13096
$ws1.Range("F4").Value = 32
$ws1.Range("F7").Value = 56
$ws1.Range("F10").Value = 13058
$ws1.Range("F11").Value = 306
$ws1.Range("F12").Value = 553
$ws1.Range("F13").Value = 8760
$ws1.Range("F14").Value = 7807
$ws1.Range("F15").Value = 214
$ws1.Range("F16").Value = 133
$ws1.Range("F18").Value = 136
$ws1.Range("F24").Value = 338

# Sheet "全部类型" updates (column F = "想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 9
$ws4.Range("F4").Value = 13097
$ws4.Range("F5").Value = 32
$ws4.Range("F8").Value = 56
$ws4.Range("F11").Value = 13058
$ws4.Range("F12").Value = 306
$ws4.Range("F13").Value = 553
$ws4.Range("F14").Value = 8760
$ws4.Range("F15").Value = 7807
$ws4.Range("F16").Value = 214
$ws4.Range("F17").Value = 133
$ws4.Range("F19").Value = 136
$ws4.Range("F27").Value = 338
